$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 689. This shifts the existing rows 689:725
# (the weekly "Vega Monumental Concepción - Plátano" records) down to
# 690:726, matching the diff's new dimension A1:T726.
$ws.Rows.Item(689).Insert()

# Populate the newly inserted row 689 with the new weekly record. The
# "constant" columns (market/product identifiers, units, origin, etc.) are
# the same values used throughout this product's block of rows.
$ws.Range("A689").Value = 11
$ws.Range("B689").Value = "Vega Monumental Concepción"
$ws.Range("C689").Value = "Bíobío"
$ws.Range("D689").Value = 45008
$ws.Range("E689").Value = 8
$ws.Range("F689").Value = "Fruta"
$ws.Range("G689").Value = 100108
$ws.Range("H689").Value = "Tropicales y subtropicales"
$ws.Range("I689").Value = 100108006
$ws.Range("J689").Value = "Plátano"
$ws.Range("K689").Value = "Sin especificar"
$ws.Range("L689").Value = "Pintón"
$ws.Range("M689").Value = 900
$ws.Range("N689").Value = 21000
$ws.Range("O689").Value = 22000
$ws.Range("P689").Value = 21556
$ws.Range("Q689").Value = "$/caja 20 kilos"
$ws.Range("R689").Value = "Ecuador"
$ws.Range("S689").Value = 1078
$ws.Range("T689").Value = 20
